$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 39970
$ws.Range("I21").Value = 39970
$ws.Range("K21").Value = 39970
$ws.Range("M21").Value = -39502
$ws.Range("H23").Value = 39970
$ws.Range("I23").Value = 39970
$ws.Range("K23").Value = 39970
$ws.Range("M23").Value = -39736
$ws.Range("H86").Value = 6582671
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 7522481.5
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 7522481.5
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -7524727.5
$ws.Range("H89").Value = 6582671
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 7522481.5
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 37612407.5
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -37623639.5
$ws.Range("H103").Value = 391.58823
$ws.Range("I103").Value = 221.41667
$ws.Range("K103").Value = 664.25001
$ws.Range("M103").Value = -78.25000999999997
$ws.Range("H111").Value = 848050.1
$ws.Range("I111").Value = 1101671.9
$ws.Range("K111").Value = 3305015.7
$ws.Range("M111").Value = -3301948.7
$ws.Range("H127").Value = 12356.083
$ws.Range("I127").Value = 12356.083
$ws.Range("K127").Value = 37068.249
$ws.Range("M127").Value = -32108.249
$ws.Range("H132").Value = 2215.8909
$ws.Range("I132").Value = 1135.3636
$ws.Range("J132").Value = 6538
$ws.Range("K132").Value = 3406.0908
$ws.Range("L132").Value = 19614
$ws.Range("M132").Value = -876.0907999999999
$ws.Range("N132").Value = -24674
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5999.5
$ws.Range("I2").Value = 5999.5
$ws.Range("K2").Value = 5999.5
$ws.Range("M2").Value = -5886.5
$ws.Range("H45").Value = 7885.0625
$ws.Range("J45").Value = 8502.833000000001
$ws.Range("L45").Value = 8502.833000000001
$ws.Range("N45").Value = -9256.833000000001
$ws.Range("H116").Value = 5999.5
$ws.Range("I116").Value = 5999.5
$ws.Range("K116").Value = 5999.5
$ws.Range("M116").Value = -3705.5
$ws.Range("H122").Value = 5536.25
$ws.Range("I122").Value = 3797.1667
$ws.Range("K122").Value = 11391.5001
$ws.Range("M122").Value = -8941.500100000001
$ws.Range("H132").Value = 5072.4653
$ws.Range("I132").Value = 1722.2142
$ws.Range("K132").Value = 5166.642599999999
$ws.Range("M132").Value = -2636.642599999999
$ws.Range("H133").Value = 1301250
$ws.Range("J133").Value = 1301250
$ws.Range("L133").Value = 1301250
$ws.Range("N133").Value = -1306310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5999.5
$ws.Range("I3").Value = 5999.5
$ws.Range("K3").Value = 5999.5
$ws.Range("M3").Value = -5885.5
$ws.Range("H22").Value = 899.5
$ws.Range("I22").Value = 899.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 899.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -726.5
$ws.Range("N22").ClearContents()
$ws.Range("H102").Value = 5543.5
$ws.Range("I102").Value = 5543.5
$ws.Range("K102").Value = 5543.5
$ws.Range("M102").Value = -2298.5
$ws.Range("H134").Value = 4811
$ws.Range("I134").Value = 3006.2307
$ws.Range("K134").Value = 9018.6921
$ws.Range("M134").Value = -6483.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5736.5713
$ws.Range("I99").Value = 4659
$ws.Range("J99").Value = 6167.6
$ws.Range("K99").Value = 4659
$ws.Range("L99").Value = 6167.6
$ws.Range("M99").Value = -3161
$ws.Range("N99").Value = -9163.6
$ws.Range("H105").Value = 1500.4
$ws.Range("I105").Value = 1444.8889
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1444.8889
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 302.1111000000001
$ws.Range("N105").Value = -5494
$ws.Range("H126").Value = 5736.5713
$ws.Range("I126").Value = 4659
$ws.Range("J126").Value = 6167.6
$ws.Range("K126").Value = 13977
$ws.Range("L126").Value = 18502.8
$ws.Range("M126").Value = -11507
$ws.Range("N126").Value = -23442.8
$ws.Range("H132").Value = 4160.9653
$ws.Range("I132").Value = 2611.1333
$ws.Range("J132").Value = 5821.5
$ws.Range("K132").Value = 7833.3999
$ws.Range("L132").Value = 17464.5
$ws.Range("M132").Value = -5303.3999
$ws.Range("N132").Value = -22524.5
$ws.Range("H134").Value = 3599.1052
$ws.Range("I134").Value = 2531.1667
$ws.Range("K134").Value = 7593.500100000001
$ws.Range("M134").Value = -5058.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 15000001
$ws.Range("J9").Value = 12500001
$ws.Range("L9").Value = 37500003
$ws.Range("N9").Value = -37500451
$ws.Range("H110").Value = 4000
$ws.Range("I110").Value = 4000
$ws.Range("K110").Value = 12000
$ws.Range("M110").Value = -7910
$ws.Range("H136").Value = 8210
$ws.Range("I136").Value = 7315
$ws.Range("K136").Value = 21945
$ws.Range("M136").Value = -16845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1337977.9
$ws.Range("I80").Value = 1004985.8
$ws.Range("K80").Value = 1004985.8
$ws.Range("M80").Value = -1003987.8
$ws.Range("H83").Value = 1337977.9
$ws.Range("I83").Value = 1004985.8
$ws.Range("K83").Value = 5024929
$ws.Range("M83").Value = -5019937
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H102").Value = 2694.5715
$ws.Range("I102").Value = 3192.6
$ws.Range("J102").Value = 1449.5
$ws.Range("K102").Value = 3192.6
$ws.Range("L102").Value = 1449.5
$ws.Range("M102").Value = -1570.6
$ws.Range("N102").Value = -4693.5
$ws.Range("H122").Value = 6218.8057
$ws.Range("I122").Value = 5388.7393
$ws.Range("J122").Value = 7687.385
$ws.Range("K122").Value = 16166.2179
$ws.Range("L122").Value = 23062.155
$ws.Range("M122").Value = -13716.2179
$ws.Range("N122").Value = -27962.155
$ws.Range("H132").Value = 326409.06
$ws.Range("I132").Value = 388276.84
$ws.Range("J132").Value = 4696.6
$ws.Range("K132").Value = 1164830.52
$ws.Range("L132").Value = 14089.8
$ws.Range("M132").Value = -1162300.52
$ws.Range("N132").Value = -19149.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1303.8572
$ws.Range("I22").Value = 1242.2858
$ws.Range("J22").Value = 1427
$ws.Range("K22").Value = 1242.2858
$ws.Range("L22").Value = 1427
$ws.Range("M22").Value = -947.2858000000001
$ws.Range("N22").Value = -2017
$ws.Range("H27").Value = 1303.8572
$ws.Range("I27").Value = 1242.2858
$ws.Range("J27").Value = 1427
$ws.Range("K27").Value = 1242.2858
$ws.Range("L27").Value = 1427
$ws.Range("M27").Value = -1135.2858
$ws.Range("N27").Value = -1641
$ws.Range("H46").Value = 3389
$ws.Range("I46").Value = 2859.524
$ws.Range("K46").Value = 2859.524
$ws.Range("M46").Value = -2671.524
$ws.Range("H82").Value = 1767.091
$ws.Range("I82").Value = 2025.8889
$ws.Range("J82").Value = 602.5
$ws.Range("K82").Value = 2025.8889
$ws.Range("L82").Value = 602.5
$ws.Range("M82").Value = -1664.8889
$ws.Range("N82").Value = -1324.5
$ws.Range("H85").Value = 1767.091
$ws.Range("I85").Value = 2025.8889
$ws.Range("J85").Value = 602.5
$ws.Range("K85").Value = 2025.8889
$ws.Range("L85").Value = 602.5
$ws.Range("M85").Value = -777.8888999999999
$ws.Range("N85").Value = -3098.5
$ws.Range("H132").Value = 5281.4
$ws.Range("I132").Value = 4174.5713
$ws.Range("K132").Value = 12523.7139
$ws.Range("M132").Value = -9993.713899999999
$ws.Range("H133").Value = 58571.43
$ws.Range("J133").Value = 58571.43
$ws.Range("L133").Value = 58571.43
$ws.Range("N133").Value = -63631.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5179
$ws.Range("I126").Value = 4965.1665
$ws.Range("K126").Value = 14895.4995
$ws.Range("M126").Value = -12425.4995
$ws.Range("H132").Value = 4112.357
$ws.Range("I132").Value = 3095.5293
$ws.Range("K132").Value = 9286.5879
$ws.Range("M132").Value = -6756.5879
$ws.Range("H140").Value = 111833.336
$ws.Range("J140").Value = 111833.336
$ws.Range("L140").Value = 111833.336
$ws.Range("N140").Value = -122193.336
